# P6 - Data Visualization - Post Udacity Reviewer and many more revisions
# Adds new "Time Spent" log entries (AB Testing / Fixing Data Visualization
# Project) around 2016-05-14/15/21, which pushes the existing row for
# 2016-05-15 onward down by one row.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Time Spent")
$ws.Activate()

# Insert a new row at 309 - everything from the old row 309 downward shifts
# down by one (dates stay the same values, just move to row+1).
$ws.Rows("309").Insert()

# Row 307 (2016-05-14): hours bumped from 8 to 9.
$ws.Cells.Item(307, 2).Value = 9

# Row 308 (2016-05-15): add a "Data Visualization Project" entry.
$ws.Cells.Item(308, 2).Value = 2
$ws.Cells.Item(308, 3).Value = "Data Visualization Project"

# Row 309 (new row, 2016-05-15): add an "AB Testing" entry.
$ws.Cells.Item(309, 1).Value = 42505
$ws.Cells.Item(309, 2).Value = 1
$ws.Cells.Item(309, 3).Value = "AB Testing"

# Row 315 (2016-05-21): add a "Fixing Data Visualization Project" entry.
$ws.Cells.Item(315, 2).Value = 6
$ws.Cells.Item(315, 3).Value = "Fixing Data Visualization Project"

# Restore the view state recorded after these edits.
$ws.Application.ActiveWindow.ScrollRow = 293
$ws.Range("B316").Select()
